$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: new auto-log entry, same config as row 25 (Q=2, D=2) ---
# Copy row 25 as a template to preserve cell types (incl. empty inlineStr cells)
$ws.Range("A25:Z25").Copy($ws.Range("A26:Z26"))

$ws.Range("A26").Value = 25
$ws.Range("T26").Value = 0.0002126383832875674
$ws.Range("U26").Value = 0.0006463956604885096
$ws.Range("V26").Value = 3.32081937789917
$ws.Range("W26").Value = 4.308806896209717
$ws.Range("X26").Value = 2.274315118789673
$ws.Range("Y26").Value = "test run with Quantum layer"
$ws.Range("Z26").Value = "AAPL, MSFT, GOOGL"

# --- Row 27: new auto-log entry, classical-only config (Q=-, D=-) like row 6 ---
# Copy row 6 as a template to preserve cell types (incl. "-" strings and empty inlineStr cells)
$ws.Range("A6:Z6").Copy($ws.Range("A27:Z27"))

$ws.Range("A27").Value = 26
$ws.Range("C27").Value = 4
$ws.Range("R27").Value = 15
$ws.Range("T27").Value = 0.0003715450874300045
$ws.Range("U27").Value = 0.007181070270041559
$ws.Range("V27").Value = 10.74910068511963
$ws.Range("W27").Value = 14.20177745819092
$ws.Range("X27").Value = 6.630932331085205
$ws.Range("Y27").Value = "test run with classical only layer"
$ws.Range("Z27").Value = "AAPL, MSFT, GOOGL"
